$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6226999999999999
$ws.Range("I2").Value = 0.8812810914468889
$ws.Range("J2").Value = 0.881281091446889
$ws.Range("M2").Value = 1.139366
$ws.Range("N2").Value = 3.418098
$ws.Range("O2").Value = 0.1546713947032042
$ws.Range("P2").Value = 0.1546713947032042
$ws.Range("Q2").Value = 0.7094832082
$ws.Range("R2").Value = 6.3853488738
$ws.Range("S2").Value = 0.1363089755396523
$ws.Range("T2").Value = 0.1363089755396523
$ws.Range("G3").Value = 0.6226999999999999
$ws.Range("I3").Value = 0.8812810914468889
$ws.Range("J3").Value = 0.881281091446889
$ws.Range("O3").Value = 0.7676983257595695
$ws.Range("P3").Value = 0.7676983257595695
$ws.Range("Q3").Value = 3.521459621766666
$ws.Range("S3").Value = 0.6765580184273428
$ws.Range("T3").Value = 0.6765580184273429
$ws.Range("G4").Value = 0.6226999999999999
$ws.Range("I4").Value = 0.8812810914468889
$ws.Range("J4").Value = 0.881281091446889
$ws.Range("M4").Value = 0.4872916666666667
$ws.Range("N4").Value = 1.461875
$ws.Range("O4").Value = 0.06615089594615092
$ws.Range("P4").Value = 0.06615089594615092
$ws.Range("Q4").Value = 0.3034365208333333
$ws.Range("R4").Value = 2.7309286875
$ws.Range("S4").Value = 0.05829753377961346
$ws.Range("T4").Value = 0.05829753377961347
$ws.Range("G5").Value = 0.6226999999999999
$ws.Range("I5").Value = 0.8812810914468889
$ws.Range("J5").Value = 0.881281091446889
$ws.Range("M5").Value = 0.08456133333333334
$ws.Range("N5").Value = 0.253684
$ws.Range("O5").Value = 0.0114793835910754
$ws.Range("P5").Value = 0.0114793835910754
$ws.Range("Q5").Value = 0.05265634226666666
$ws.Range("R5").Value = 0.4739070804
$ws.Range("S5").Value = 0.01011656370028044
$ws.Range("T5").Value = 0.01011656370028044
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04989433333333334
$ws.Range("H6").Value = 0.149683
$ws.Range("I6").Value = 0.07061334918422178
$ws.Range("J6").Value = 0.07061334918422178
$ws.Range("M6").Value = 1.139366
$ws.Range("N6").Value = 3.418098
$ws.Range("O6").Value = 0.1546713947032042
$ws.Range("P6").Value = 0.1546713947032042
$ws.Range("Q6").Value = 0.05684790699266668
$ws.Range("R6").Value = 0.511631162934
$ws.Range("S6").Value = 0.01092186520298795
$ws.Range("T6").Value = 0.01092186520298795
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04989433333333334
$ws.Range("H7").Value = 0.149683
$ws.Range("I7").Value = 0.07061334918422178
$ws.Range("J7").Value = 0.07061334918422178
$ws.Range("O7").Value = 0.7676983257595695
$ws.Range("P7").Value = 0.7676983257595695
$ws.Range("Q7").Value = 0.2821597562041112
$ws.Range("R7").Value = 2.539437805837
$ws.Range("S7").Value = 0.05420974994500292
$ws.Range("T7").Value = 0.05420974994500292
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.04989433333333334
$ws.Range("H8").Value = 0.149683
$ws.Range("I8").Value = 0.07061334918422178
$ws.Range("J8").Value = 0.07061334918422178
$ws.Range("M8").Value = 0.4872916666666667
$ws.Range("N8").Value = 1.461875
$ws.Range("O8").Value = 0.06615089594615092
$ws.Range("P8").Value = 0.06615089594615092
$ws.Range("Q8").Value = 0.02431309284722222
$ws.Range("R8").Value = 0.218817835625
$ws.Range("S8").Value = 0.004671136314294676
$ws.Range("T8").Value = 0.004671136314294676
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.04989433333333334
$ws.Range("H9").Value = 0.149683
$ws.Range("I9").Value = 0.07061334918422178
$ws.Range("J9").Value = 0.07061334918422178
$ws.Range("M9").Value = 0.08456133333333334
$ws.Range("N9").Value = 0.253684
$ws.Range("O9").Value = 0.0114793835910754
$ws.Range("P9").Value = 0.0114793835910754
$ws.Range("Q9").Value = 0.004219131352444445
$ws.Range("R9").Value = 0.037972182172
$ws.Range("S9").Value = 0.0008105977219362329
$ws.Range("T9").Value = 0.0008105977219362332
$ws.Range("G10").Value = 0.03399066666666666
$ws.Range("I10").Value = 0.04810555936888933
$ws.Range("J10").Value = 0.04810555936888933
$ws.Range("M10").Value = 1.139366
$ws.Range("N10").Value = 3.418098
$ws.Range("O10").Value = 0.1546713947032042
$ws.Range("P10").Value = 0.1546713947032042
$ws.Range("Q10").Value = 0.03872780991733333
$ws.Range("R10").Value = 0.348550289256
$ws.Range("S10").Value = 0.007440553960563902
$ws.Range("T10").Value = 0.007440553960563903
$ws.Range("G11").Value = 0.03399066666666666
$ws.Range("I11").Value = 0.04810555936888933
$ws.Range("J11").Value = 0.04810555936888933
$ws.Range("O11").Value = 0.7676983257595695
$ws.Range("P11").Value = 0.7676983257595695
$ws.Range("Q11").Value = 0.1922221939675555
$ws.Range("S11").Value = 0.03693055738722391
$ws.Range("T11").Value = 0.03693055738722392
$ws.Range("G12").Value = 0.03399066666666666
$ws.Range("I12").Value = 0.04810555936888933
$ws.Range("J12").Value = 0.04810555936888933
$ws.Range("M12").Value = 0.4872916666666667
$ws.Range("N12").Value = 1.461875
$ws.Range("O12").Value = 0.06615089594615092
$ws.Range("P12").Value = 0.06615089594615092
$ws.Range("Q12").Value = 0.01656336861111111
$ws.Range("R12").Value = 0.1490703175
$ws.Range("S12").Value = 0.003182225852242784
$ws.Range("T12").Value = 0.003182225852242784
$ws.Range("G13").Value = 0.03399066666666666
$ws.Range("I13").Value = 0.04810555936888933
$ws.Range("J13").Value = 0.04810555936888933
$ws.Range("M13").Value = 0.08456133333333334
$ws.Range("N13").Value = 0.253684
$ws.Range("O13").Value = 0.0114793835910754
$ws.Range("P13").Value = 0.0114793835910754
$ws.Range("Q13").Value = 0.002874296094222222
$ws.Range("R13").Value = 0.025868664848
$ws.Range("S13").Value = 0.0005522221688587316
$ws.Range("T13").Value = 0.0005522221688587318
